$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells for new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the style of an existing header cell (e.g. AC1) to the new headers
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Re-set values since paste special with formats only shouldn't overwrite them, but ensure values remain correct
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the Wins/Losses/Ties data for rows 2 through 40
for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 30).Value = 78   # AD = column 30
    $ws.Cells.Item($r, 31).Value = 84   # AE = column 31
    $ws.Cells.Item($r, 32).Value = 0    # AF = column 32
}
